$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.727.03"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.592.33"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.94"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.30"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.818.43"
$ws.Range("E12").Value = "  -0.37%  "
$ws.Range("D13").Value = "1.582.91"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  -0.79%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "27.716.19"
$ws.Range("E16").Value = "  +0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.26"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "218.20"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").Value = "0.0₃0695"
$ws.Range("E19").Value = "  +0.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.35"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.74"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.81"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.05"
$ws.Range("E26").Value = "  +5.13%  "
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.11"
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +1.60%  "
$ws.Range("E32").Value = "  -2.48%  "
$ws.Range("D33").Value = "1.382.83"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.34"
$ws.Range("E37").Value = "  +0.80%  "
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.535"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("E40").Value = "  +1.32%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.985"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.47"
$ws.Range("E43").Value = "  +0.58%  "
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "1.730.58"
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.09"
$ws.Range("E48").Value = "  -1.86%  "
$ws.Range("E49").Value = "  +0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("E51").Value = "  -0.18%  "
